$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.746.44"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "2.648.43"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'112.83"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").Value = "'327.98"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").Value = "'0.524"
$ws.Range("E7").Value = "  -0.95%  "

$ws.Range("D8").Value = "'1.00"

$ws.Range("D9").Value = "'0.551"
$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").Value = "'39.81"
$ws.Range("E10").Value = "  -3.08%  "

$ws.Range("D11").Value = "'20.00"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").Value = "'7.58"
$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("D15").Value = "3.066.71"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "2.654.53"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "'0.863"

$ws.Range("D18").Value = "49.761.77"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").Value = "'13.43"
$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").Value = "'6.70"
$ws.Range("E21").Value = "  -0.96%  "

$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  -0.97%  "

$ws.Range("D23").Value = "'269.17"
$ws.Range("E23").Value = "  -2.25%  "

$ws.Range("D24").Value = "'69.28"
$ws.Range("E24").Value = "  -4.37%  "

$ws.Range("D25").Value = "'2.56"
$ws.Range("E25").Value = "  -0.98%  "

$ws.Range("D26").Value = "'26.19"
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  +1.80%  "

$ws.Range("E29").Value = "  -0.87%  "

$ws.Range("D30").Value = "'0.139"
$ws.Range("E30").Value = "  -1.94%  "

$ws.Range("D31").Value = "'35.07"
$ws.Range("E31").Value = "  -4.34%  "

$ws.Range("D32").Value = "'49.59"
$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("D33").Value = "'5.51"
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("E34").Value = "  +0.54%  "

$ws.Range("D35").Value = "'19.17"
$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").Value = "'4.96"
$ws.Range("E37").Value = "  -1.97%  "

$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("D39").Value = "'3.14"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'23.66"
$ws.Range("E40").Value = "  +6.79%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'128.38"
$ws.Range("E41").Value = "  +3.29%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'2.30"
$ws.Range("E42").Value = "  +3.61%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0344"
$ws.Range("E43").Value = "  +7.96%  "

$ws.Range("E44").Value = "  -0.60%  "

$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "2.060.66"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("E47").Value = "  +7.53%  "

$ws.Range("E48").Value = "  -3.01%  "

$ws.Range("D49").Value = "'8.96"
$ws.Range("E49").Value = "  -2.17%  "

$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("D51").Value = "'59.17"
$ws.Range("E51").Value = "  -2.95%  "

